$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo in existing row 4 label: "Maquina #D" -> "Maquina 3D"
$ws.Range("A4").Value = "Maquina 3D"

# Widen column A a bit (matches new <cols> override)
$ws.Columns.Item(1).ColumnWidth = 16.8

# New machine rows appended after the existing data (rows 5-12)
$data = @(
    @("Escaner", 1234),
    @("Multifuncional Xerox", 5678),
    @("Multifuncional Canon", 91011),
    @("Multifuncional HP", 1213),
    @("Multifuncional Epson", 1415),
    @("Multifuncional Ricoh", 1617),
    @("Multifuncional Brother", 1819),
    @("Multifuncional Samsung", 2021)
)

$row = 5
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $row = $row + 1
}

# Column A (rows 5-12) reuses the header's style (s=1)
$ws.Range("A1").Copy()
$ws.Range("A5:A12").PasteSpecial(-4122)

# B5 reuses the data-body style (s=2), like A2:B4 originally had
$ws.Range("B2").Copy()
$ws.Range("B5").PasteSpecial(-4122)

# B6:B12 reuse the header-row style (s=1), matching the source workbook
$ws.Range("A1").Copy()
$ws.Range("B6:B12").PasteSpecial(-4122)
